$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 7.765487334488884
$ws.Range("D2").Value2 = 7.729717733732095
$ws.Range("E2").Value2 = 12.72596180790718
$ws.Range("F2").Value2 = 44.58921583010547
$ws.Range("G2").Value2 = 57.8414135015265
$ws.Range("H2").Value2 = 19.80072001141792
$ws.Range("J2").Value2 = 10.109041464798
$ws.Range("N2").Value2 = 18.20530214749979

$ws.Range("B3").Value2 = 7.690488708080738
$ws.Range("D3").Value2 = 7.731607383498005
$ws.Range("E3").Value2 = 12.7347732422308
$ws.Range("F3").Value2 = 43.77798400427347
$ws.Range("G3").Value2 = 56.1295405628279
$ws.Range("H3").Value2 = 19.59404963810021
$ws.Range("J3").Value2 = 10.11490316061085
$ws.Range("N3").Value2 = 17.93039231733595

$ws.Range("B4").Value2 = 7.645832371061139
$ws.Range("D4").Value2 = 7.732971663349161
$ws.Range("E4").Value2 = 12.74289888234467
$ws.Range("F4").Value2 = 43.2888058150247
$ws.Range("G4").Value2 = 55.074807882116
$ws.Range("H4").Value2 = 19.47262906352406
$ws.Range("J4").Value2 = 10.12069478094539
$ws.Range("N4").Value2 = 17.76180298315366

$ws.Range("B5").Value2 = 7.628005955137633
$ws.Range("D5").Value2 = 7.733578766040622
$ws.Range("E5").Value2 = 12.74689087757529
$ws.Range("F5").Value2 = 43.09197246160743
$ws.Range("G5").Value2 = 54.64481053314584
$ws.Range("H5").Value2 = 19.42457083818628
$ws.Range("J5").Value2 = 10.12360443381637
$ws.Range("N5").Value2 = 17.69323820308988

$ws.Range("B6").Value2 = 7.625068927627218
$ws.Range("D6").Value2 = 7.733682659711756
$ws.Range("E6").Value2 = 12.74759479490273
$ws.Range("F6").Value2 = 43.05944780799597
$ws.Range("G6").Value2 = 54.57342058210207
$ws.Range("H6").Value2 = 19.4166779099519
$ws.Range("J6").Value2 = 10.12412071270919
$ws.Range("N6").Value2 = 17.681863828119

$ws.Range("B7").Value2 = 7.645590426591721
$ws.Range("D7").Value2 = 7.732979644053993
$ws.Range("E7").Value2 = 12.74294996641711
$ws.Range("F7").Value2 = 43.2861407353301
$ws.Range("G7").Value2 = 55.06900854887076
$ws.Range("H7").Value2 = 19.47197512071688
$ws.Range("J7").Value2 = 10.12073179906159
$ws.Range("N7").Value2 = 17.76087762816193

$ws.Range("B8").Value2 = 7.739351513268095
$ws.Range("D8").Value2 = 7.730326842399959
$ws.Range("E8").Value2 = 12.72843522590959
$ws.Range("F8").Value2 = 44.30781285096198
$ws.Range("G8").Value2 = 57.25234060176211
$ws.Range("H8").Value2 = 19.72835136274256
$ws.Range("J8").Value2 = 10.11060648771349
$ws.Range("N8").Value2 = 18.11052099476146

$ws.Range("B9").Value2 = 7.9333151054759
$ws.Range("D9").Value2 = 7.726751176464113
$ws.Range("E9").Value2 = 12.72160251247386
$ws.Range("F9").Value2 = 46.37024980785385
$ws.Range("G9").Value2 = 61.47491567459172
$ws.Range("H9").Value2 = 20.27243875824717
$ws.Range("J9").Value2 = 10.10822160243079
$ws.Range("N9").Value2 = 18.7941399855402

$ws.Range("B10").Value2 = 8.080666255167932
$ws.Range("D10").Value2 = 7.725127013220942
$ws.Range("E10").Value2 = 12.72987491738833
$ws.Range("F10").Value2 = 47.90568748517494
$ws.Range("G10").Value2 = 64.50414412291448
$ws.Range("H10").Value2 = 20.69436942026452
$ws.Range("J10").Value2 = 10.11721207091131
$ws.Range("N10").Value2 = 19.29036450727469

$ws.Range("B11").Value2 = 8.148469166187692
$ws.Range("D11").Value2 = 7.724608544374606
$ws.Range("E11").Value2 = 12.73654248767226
$ws.Range("F11").Value2 = 48.6053310187205
$ws.Range("G11").Value2 = 65.85944379355729
$ws.Range("H11").Value2 = 20.89045543231269
$ws.Range("J11").Value2 = 10.12365035974141
$ws.Range("N11").Value2 = 19.51383106981899

$ws.Range("B12").Value2 = 8.17423166240809
$ws.Range("D12").Value2 = 7.724444128018039
$ws.Range("E12").Value2 = 12.73948595990004
$ws.Range("F12").Value2 = 48.87017583043146
$ws.Range("G12").Value2 = 66.368886747877
$ws.Range("H12").Value2 = 20.96524579847746
$ws.Range("J12").Value2 = 10.12642695060378
$ws.Range("N12").Value2 = 19.59804885977593

$ws.Range("B13").Value2 = 8.168679758780495
$ws.Range("D13").Value2 = 7.724478114804402
$ws.Range("E13").Value2 = 12.73883340025739
$ws.Range("F13").Value2 = 48.81314504986548
$ws.Range("G13").Value2 = 66.25934484517246
$ws.Range("H13").Value2 = 20.9491154314576
$ws.Range("J13").Value2 = 10.12581389385634
$ws.Range("N13").Value2 = 19.57993023154885

$ws.Range("B14").Value2 = 8.150587055476455
$ws.Range("D14").Value2 = 7.724594376744808
$ws.Range("E14").Value2 = 12.73677625479573
$ws.Range("F14").Value2 = 48.62712332969024
$ws.Range("G14").Value2 = 65.90143448747506
$ws.Range("H14").Value2 = 20.89659806435153
$ws.Range("J14").Value2 = 10.12387200243344
$ws.Range("N14").Value2 = 19.52076825227035

$ws.Range("B15").Value2 = 8.139515373066576
$ws.Range("D15").Value2 = 7.724669753947895
$ws.Range("E15").Value2 = 12.73557073299515
$ws.Range("F15").Value2 = 48.51315967612338
$ws.Range("G15").Value2 = 65.68169744160305
$ws.Range("H15").Value2 = 20.86449779800469
$ws.Range("J15").Value2 = 10.12272664843227
$ws.Range("N15").Value2 = 19.48447494387074

$ws.Range("B16").Value2 = 8.076248655077107
$ws.Range("D16").Value2 = 7.72516535002237
$ws.Range("E16").Value2 = 12.72949770781503
$ws.Range("F16").Value2 = 47.85996489177818
$ws.Range("G16").Value2 = 64.41507003753819
$ws.Range("H16").Value2 = 20.68163311876464
$ws.Range("J16").Value2 = 10.11683864795549
$ws.Range("N16").Value2 = 19.27570805562912

$ws.Range("B17").Value2 = 8.037617191288176
$ws.Range("D17").Value2 = 7.725525983146547
$ws.Range("E17").Value2 = 12.72651675488256
$ws.Range("F17").Value2 = 47.45935932342474
$ws.Range("G17").Value2 = 63.63183684729849
$ws.Range("H17").Value2 = 20.5704716703448
$ws.Range("J17").Value2 = 10.11382872043847
$ws.Range("N17").Value2 = 19.14699881872645

$ws.Range("B18").Value2 = 8.015471503473739
$ws.Range("D18").Value2 = 7.725754137937617
$ws.Range("E18").Value2 = 12.72507552557142
$ws.Range("F18").Value2 = 47.22906635946156
$ws.Range("G18").Value2 = 63.17923164279156
$ws.Range("H18").Value2 = 20.50692855796396
$ws.Range("J18").Value2 = 10.11231850326613
$ws.Range("N18").Value2 = 19.07275961144547

$ws.Range("B19").Value2 = 8.00798681698603
$ws.Range("D19").Value2 = 7.725834940556858
$ws.Range("E19").Value2 = 12.72463446245567
$ws.Range("F19").Value2 = 47.15112272176361
$ws.Range("G19").Value2 = 63.02564133938778
$ws.Range("H19").Value2 = 20.48548342974222
$ws.Range("J19").Value2 = 10.11184509972871
$ws.Range("N19").Value2 = 19.04759001474568

$ws.Range("B20").Value2 = 8.041722081219515
$ws.Range("D20").Value2 = 7.725485446234032
$ws.Range("E20").Value2 = 12.72680578544053
$ws.Range("F20").Value2 = 47.50199349396804
$ws.Range("G20").Value2 = 63.71543563327825
$ws.Range("H20").Value2 = 20.58226464743871
$ws.Range("J20").Value2 = 10.11412625151573
$ws.Range("N20").Value2 = 19.16072234870301

$ws.Range("B21").Value2 = 8.155899150234964
$ws.Range("D21").Value2 = 7.724559359728069
$ws.Range("E21").Value2 = 12.7373691210164
$ws.Range("F21").Value2 = 48.6817670345036
$ws.Range("G21").Value2 = 66.00666786147329
$ws.Range("H21").Value2 = 20.91200959515086
$ws.Range("J21").Value2 = 10.12443318921646
$ws.Range("N21").Value2 = 19.53815713259445

$ws.Range("B22").Value2 = 8.231016616546755
$ws.Range("D22").Value2 = 7.724140255687939
$ws.Range("E22").Value2 = 12.74671310399783
$ws.Range("F22").Value2 = 49.45214975722511
$ws.Range("G22").Value2 = 67.48188650198634
$ws.Range("H22").Value2 = 21.13061830241064
$ws.Range("J22").Value2 = 10.13314293819171
$ws.Range("N22").Value2 = 19.78244317294328

$ws.Range("B23").Value2 = 8.19088749503514
$ws.Range("D23").Value2 = 7.724346829107239
$ws.Range("E23").Value2 = 12.74150251249817
$ws.Range("F23").Value2 = 49.0411268991405
$ws.Range("G23").Value2 = 66.69672593583009
$ws.Range("H23").Value2 = 21.01367872353836
$ws.Range("J23").Value2 = 10.12831357741493
$ws.Range("N23").Value2 = 19.65230661566007

$ws.Range("B24").Value2 = 8.039866059139435
$ws.Range("D24").Value2 = 7.725503708133549
$ws.Range("E24").Value2 = 12.72667426583628
$ws.Range("F24").Value2 = 47.48271850152633
$ws.Range("G24").Value2 = 63.67764779024028
$ws.Range("H24").Value2 = 20.57693190082308
$ws.Range("J24").Value2 = 10.11399105181379
$ws.Range("N24").Value2 = 19.15451869034863

$ws.Range("B25").Value2 = 7.879895846024556
$ws.Range("D25").Value2 = 7.727543224494022
$ws.Range("E25").Value2 = 12.72112336097148
$ws.Range("F25").Value2 = 45.80763838381693
$ws.Range("G25").Value2 = 60.34292603907729
$ws.Range("H25").Value2 = 20.12112573215001
$ws.Range("J25").Value2 = 10.10698609646884
$ws.Range("N25").Value2 = 18.60992361498231
